$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates derived from the authoritative diff.
# Values that parse as plain numbers need NumberFormat forced to Text
# first, otherwise Excel auto-converts the assigned string into a
# numeric cell (these columns store numbers formatted with '.' as a
# thousands separator, e.g. '70.897.67', so they must stay text).

$ws.Range("D2").Value = "70.897.67"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "3.573.86"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.98"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.48"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("E7").Value = "  +2.22%  "
$ws.Range("D8").Value = "3.562.70"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("E10").Value = "  +16.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.651"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.76"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").Value = "  +5.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.57"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "4.141.06"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "70.917.14"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.601.07"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.24"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.80"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "582.49"
$ws.Range("E20").Value = "  +6.70%  "
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.72"
$ws.Range("E23").Value = "  -5.31%  "
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.13"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.22"
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.98"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.77"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.34"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.81"
$ws.Range("E34").Value = "  +22.19%  "
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("E36").Value = "  +8.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "538.13"
$ws.Range("E37").Value = "  -3.64%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0827"
$ws.Range("E38").Value = "  +7.85%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.412"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.17"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "3.621.00"
$ws.Range("E42").Value = "  +10.28%  "
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.46"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0476"
$ws.Range("E45").Value = "  +7.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.50"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.41"
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.46"
$ws.Range("E51").Value = "  +6.09%  "
